$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$valLOT2052 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"
$valLOT2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

# Swap the contents of row 23 (currently LOT2052) and row 24 (currently LOT2028)
$ws.Range("B23").Value = $valLOT2028
$ws.Range("C23").Value = $valLOT2028
$ws.Range("B24").Value = $valLOT2052
$ws.Range("C24").Value = $valLOT2052
